# chore: minor edits
# - Recolor the "Recall@N" (D) column fill from yellow to red
# - Clear out the "Novelty@N" column (E) header + values (column kept, just emptied)
# - Drop the now-unused "Novelty@N" hyperlink, keep the "Recall@N" one on D1
# - Update the active selection to H19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column D (Recall@N values, rows 2-49): fill red instead of yellow.
#    255 = RGB(255,0,0) in Excel's BGR-packed Long color value.
$ws.Range("D2:D49").Interior.Color = 255

# 2) Column E ("Novelty@N"): clear the header cell and all data values,
#    but keep the cells (and their existing number/fill style) in place.
$ws.Range("E1").ClearContents() | Out-Null
$ws.Range("E2:E49").ClearContents() | Out-Null

# 3) Hyperlinks: the sheet only exposes collection-level Delete() reliably,
#    so drop every hyperlink and re-create the one that should survive
#    (D1 -> mailto:Recall@N). This also removes the stale E1 -> Novelty@N link.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("D1"), "mailto:Recall@N") | Out-Null
# Re-adding the hyperlink reapplies Excel's built-in "Hyperlink" cell style;
# put D1 back on its original style so no spurious style entry is created.
$ws.Range("D1").Style = "Hyperlink"

# 4) Move the active selection from G19 to H19.
$ws.Range("H19").Select() | Out-Null
